$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Column X = "Rodada 23", Column Y = "Rodada 24" -- newly reported round results
# for each team (rows 2-21).
$values = @{
    2  = @(99.2,    39.58)
    3  = @(91.1,    39.06)
    4  = @(65.02,   80.86)
    5  = @(98,      71.28)
    6  = @(91.8,    62.84)
    7  = @(100.47,  43.63)
    8  = @(106.5,   85.18)
    9  = @(69.95,   49.81)
    10 = @(66.45,   69.45)
    11 = @(110.47,  45.83)
    12 = @(31.14,   40.6)
    13 = @(93,      19.41)
    14 = @(61.6,    46.11)
    15 = @(23.69,   41.54)
    16 = @(78.9,    63.31)
    17 = @(91.05,   71.58)
    18 = @(67.23,   58.61)
    19 = @(59.45,   43.35)
    20 = @(59.07,   49.08)
    21 = @(84.54,   46.54)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 24).Value = $pair[0]
    $ws.Cells.Item($row, 25).Value = $pair[1]
}

# Update the active selection left on the sheet after the edits.
$ws.Activate()
$ws.Range("AB15").Select()
